$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the delta CD inputs (O4:O6)
$ws.Range("O4").Value = 0.03
$ws.Range("O5").Value = 0.08
$ws.Range("O6").Value = 0.03

# Update the new Cl max (H7:J9) inputs
$ws.Range("H7").Value = 2.7
$ws.Range("I7").Value = 1.9
$ws.Range("J7").Value = 1.3

$ws.Range("H8").Value = 3
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 2.5
$ws.Range("J9").Value = 1.7

# Recalculate the whole workbook so all dependent formulas and chart caches refresh
$excel.CalculateFullRebuild()

# Update the saved view state for Sheet1
$ws.Activate()
$av = $excel.ActiveWindow
$av.ScrollColumn = 1
$av.ScrollRow = 10
$av.Zoom = 130
$ws.Range("E69").Select()
